$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-12 03:19:39"

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-12 03:19:44"

# Overview mirrors each locale's status (E=zh-cn, F=de-de) and the
# latest handoff xliff generation timestamp (G).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-12 03:19:44"

# Widen the status/locale columns to fit the new, longer status text
# ("Ready for handoff" vs "In Translation"), matching the width Excel
# settles on for these linked columns across the three sheets.
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
